# Add a new "query (45)" row to the Internas workbook:
#  - rename the worksheet and the query table / defined name from
#    "query (40)" generation 40 to generation 45
#  - extend the ListObject (Tabela_query__40 -> Tabela_query__45) over the
#    new row
#  - append one new record (row 15) with the data supplied by Power Query

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename sheet, query table and defined name: "query" -> "query (45)"
# ---------------------------------------------------------------------
$ws.Name = "query (45)"

$nm = $wb.Names.Item(1)
$nm.Name = "query__45"
$nm.RefersTo = "='query (45)'!`$A`$1:`$H`$15"

# ---------------------------------------------------------------------
# 2. Append the new row of data (row 15)
# ---------------------------------------------------------------------
$newText = "Cadastro de ação a pedido de: Pedro Igor Grilo de Oliveira Carvalho`nE-mail - AÇÃO OFERTA CASHBACK | INAUGURAÇÃO BAP GRAVATINHA`nAbaixo estão as informações da ação:`nPrecisamos de apoio para colocarmos uma ação para o posto abaixo para amanhã:`n101785 - CENTRO AUTOMOTIVO BAP GRAVATINHA LTDA`nCNPJ: 61620792000191 | SP - Santo Andre`nJd Bela Vista`n700 pontos por 8% de cashback em todos os combustíveis."

# Column A ("Título") stays blank for every record in this list, same as
# the existing rows - just give it the same text-number-format as the rest
# of the column.
$ws.Range("A15").NumberFormat = "@"

$ws.Range("B15").Value = $newText
$ws.Range("B15").WrapText = $true

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "Cadastro de Ações"

$ws.Range("D15").NumberFormat = "m/d/yy h:mm"
$ws.Range("D15").Value = 46009.631944444445

$ws.Range("E15").NumberFormat = "m/d/yy h:mm"
$ws.Range("E15").Value = 46009.635416666664

$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "Larissa"

$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "Item"

$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Consolida   Demandas Internas"

# Match the row height Excel computed for the wrapped paragraph text.
$ws.Rows(15).RowHeight = 304.5

# ---------------------------------------------------------------------
# 3. Grow the table (ListObject) so it covers the new row and rename it
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H15"))
$lo.Name = "Tabela_query__45"

# ---------------------------------------------------------------------
# 4. Put the selection back on A1 (matches the saved file - no stale
#    B3 selection left over from editing)
# ---------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
